$d = $word.ActiveDocument

# Remove the first of the two trailing empty "single-underline" paragraphs
$null = $d.Paragraphs(28).Range.Delete()

# Paragraph 27: "Git Push -u origin master  " -> spell-checked run split
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:lastRenderedPageBreak/><w:t>Git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Push</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> -u </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>origin</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> master  </w:t></w:r></w:p>'
$null = $d.Paragraphs(27).Range.InsertXML($xml)

# Paragraph 26 (empty) -> empty + "Git branch -m main" + "Git branch -m master" (split)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>branch</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> -m </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>main</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>branch</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> -m master</w:t></w:r></w:p>'
$null = $d.Paragraphs(26).Range.InsertXML($xml)

# Paragraph 25: "Git branch -m master" -> empty + "Receber mudancas..." + "Git pull" (split)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Receber mudanças que foram feitas no repositório:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pull</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$null = $d.Paragraphs(25).Range.InsertXML($xml)

# Paragraph 24: "Git branch -m main" -> "Enviar arquivos..." + "Git Push " (split)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Enviar arquivos para servidor do </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Push</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>'
$null = $d.Paragraphs(24).Range.InsertXML($xml)

# Paragraph 22: "Git commit -a -m..." -> spell-checked run split (Git/commit), trailing runs kept
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>commit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">-a </w:t></w:r><w:r><w:t>-m “Estou enviando todos os arquivos”</w:t></w:r></w:p>'
$null = $d.Paragraphs(22).Range.InsertXML($xml)

# Paragraph 19: "Git  commit ..." -> spell-checked run split (Git/commit), trailing runs kept
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">  </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>commit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> “nome do arquivo”</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve"> -</w:t></w:r><w:r><w:t>m “</w:t></w:r><w:r><w:t>Estou enviando somente nome do arquivo</w:t></w:r><w:r><w:t>”</w:t></w:r></w:p>'
$null = $d.Paragraphs(19).Range.InsertXML($xml)

# Paragraph 16: "Git add ." -> spell-checked run split
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>add</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> .</w:t></w:r></w:p>'
$null = $d.Paragraphs(16).Range.InsertXML($xml)

# Paragraph 13: Git add "..." -> spell-checked run split
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>add</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> “nome do arquivo”</w:t></w:r></w:p>'
$null = $d.Paragraphs(13).Range.InsertXML($xml)

# Paragraph 10: "Git init" -> spell-checked run split
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>init</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$null = $d.Paragraphs(10).Range.InsertXML($xml)

# Paragraph 7: "Git status" -> spell-checked run split
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> status</w:t></w:r></w:p>'
$null = $d.Paragraphs(7).Range.InsertXML($xml)

# Paragraph 4: "Git -version" -> spell-checked run split
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> –</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>version</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$null = $d.Paragraphs(4).Range.InsertXML($xml)

# Paragraph 1: "Comandos Fundamentais do GIt" -> spell-checked run split
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Comandos Fundamentais do </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>GIt</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$null = $d.Paragraphs(1).Range.InsertXML($xml)
